$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table for rows 2..17 (A=index, B=name, C=from_bus, D=to_bus, E=in_service)
# "line7" and "line8" are newly inserted before the "extr*" entries, which
# pushes every line/extr label from row 8 onward down by two rows.
$names = @("line1","line2","line3","line4","line5","line6","line7","line8","extr1","extr2","extr3","extr4","extr5","extr6","extr7","extr8")
$cVals = @(7, 9, 8, 8, 10, 12, 14, 16, 5, 5, 10, 7, 9, 7, 5, 8)
$dVals = @(9, 8, 10, 11, 5, 8, 11, 9, 12, 9, 11, 8, 11, 11, 7, 5)
$eVals = @($true, $false, $true, $true, $true, $true, $true, $true, $true, $true, $true, $false, $false, $true, $false, $false)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}

# Rows 16 and 17 are brand new -- give column A the same formatting
# (bold / bordered / centered style) that the rest of column A already uses.
$ws.Range("A2").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)

# Restore the values (PasteSpecial of formats only shouldn't disturb them,
# but make sure they are correct regardless).
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(17, 1).Value = 15
